{"js": "// Replace the entire body content with the new iPhone 14 Pro Max\n// specifications text. The target document is a single paragraph\n// containing one run whose content is a flat sequence of <w:t> text\n// nodes separated by <w:br/> line breaks (manual line breaks, not\n// paragraph marks). We rebuild that exact structure via insertOoxml\n// so the run/break layout matches exactly (rather than using\n// insertText/insertBreak calls, which would split the content across\n// multiple runs).\n\n  const segments = [\n    { t: \"# iPhone 14 Pro Max Specifications\" },\n    { br: true },\n    { br: true },\n    { t: \"## Display:\" },\n    { br: true },\n    { t: \"- 6.7-inch Super Retina XDR display\" },\n    { br: true },\n    { t: \"- Always-On and ProMotion technology with adaptive 120Hz refresh rate\" },\n    { br: true },\n    { t: \"- Oleophobic coating\" },\n    { br: true },\n    { t: \"- 2796x1290 resolution at 460 ppi\" },\n    { br: true },\n    { t: \"- HDR support\" },\n    { br: true },\n    { t: \"## \" },\n    { br: true },\n    { t: \"## Camera:\" },\n    { br: true },\n    { t: \"- 48MP Main camera with 24mm \u0192/1.78 aperture\" },\n    { br: true },\n    { t: \"- 12MP Ultra Wide 13mm \u0192/2.2 aperture camera\" },\n    { br: true },\n    { t: \"- 12MP 2x Telephoto 48mm \u0192/1.78 aperture camera enabled by quad-pixel sensor\" },\n    { br: true },\n    { t: \"- 12MP 3x Telephoto 77mm \u0192/2.8 aperture camera\" },\n    { br: true },\n    { t: \"- Night mode portraits enabled by LiDAR Scanner\" },\n    { br: true },\n    { t: \"- Portrait mode with advanced bokeh and Depth Control\" },\n    { br: true },\n    { t: \"- 4K video recording at 24, 25, 30 or 60 fps\" },\n    { br: true },\n    { t: \"## \" },\n    { br: true },\n    { t: \"## Body:\" },\n    { br: true },\n    { t: \"- Ceramic Shield front, textured matte glass back and stainless steel design\" },\n    { br: true },\n    { t: \"- Resists splashes, water and dust with IP68 rating\" },\n    { br: true },\n    { t: \"- Available in Space Black, Silver, Gold and Deep Purple\" },\n    { br: true },\n    { t: \"- Dimensions: 3.05 inches x 6.33 inches x 0.31 inch (77.6 mm x 160.7 mm x 7.85 mm)\" },\n    { br: true },\n    { t: \"- Weight: 8.47 ounces (240 grams)\" },\n    { br: true },\n    { t: \"## \" },\n    { br: true },\n    { t: \"## Power:\" },\n    { br: true },\n    { t: \"- Built-in rechargeable lithium-ion battery with up to 29 hours video playback\" },\n    { br: true },\n    { t: \"- MagSafe and Qi wireless charging up to 15W\" },\n    { br: true },\n    { t: \"- Fast charging capable: up to 50% in 30 minutes with 20W adapter or higher\" },\n    { br: true },\n    { t: \"## \" },\n    { br: true },\n    { t: \"## Storage:\" },\n    { br: true },\n    { t: \"- 128GB, 256GB, 512GB and 1TB versions available\" },\n    { br: true },\n    { t: \"## \" },\n    { br: true },\n    { t: \"## Connectivity:\" },\n    { br: true },\n    { t: \"- 5G, 4G LTE, Wi-Fi 6, Bluetooth\u00ae and NFC\" },\n    { br: true },\n    { t: \"- Dual eSIM support, no physical SIM slot\" },\n    { br: true },\n    { t: \"- GPS, GLONASS, Galileo, QZSS and BeiDou\" },\n    { br: true },\n    { t: \"## \" },\n    { br: true },\n    { t: \"## Price:\" },\n    { br: true },\n    { t: \"- $1,099 (128GB), $1,199 (256GB), $1,399 (512GB), $1,599 (1TB)\" },\n  ];\n\n// XML-escape helper for text inserted into <w:t> elements.\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build the flat run content: <w:t>/<w:br/> children, preserving\n// leading/trailing spaces with xml:space=\"preserve\".\nlet runInner = \"\";\nfor (const seg of segments) {\n  if (seg.br) {\n    runInner += \"<w:br/>\";\n  } else {\n    const text = seg.t;\n    const esc = xmlEscape(text);\n    const preserve = text !== text.trim() || text.length === 0;\n    runInner += preserve\n      ? `<w:t xml:space=\"preserve\">${esc}</w:t>`\n      : `<w:t>${esc}</w:t>`;\n  }\n}\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>${runInner}</w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst body = context.document.body;\nconst wholeRange = body.getRange(Word.RangeLocation.whole);\nwholeRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the entire document content with the new iPhone 14 Pro Max\n# specifications text. The target document is a single paragraph\n# containing one run whose content is a flat sequence of <w:t> text\n# nodes separated by <w:br/> manual line breaks (not paragraph marks).\n# We rebuild that exact structure by calling Range.InsertXML() with a\n# flat-OPC WordprocessingML fragment on the whole-document Content\n# range -- InsertXML REPLACES that range's contents in place, and\n# (unlike assigning Range.Text) only adds xml:space=\"preserve\" to the\n# <w:t> elements that actually need it.\n\n$d = $word.ActiveDocument\n\n# Each entry is either a run of text (Text) or a manual line break (Br).\n$segments = @(\n  @{ Text = '# iPhone 14 Pro Max Specifications' },\n  @{ Br = $true },\n  @{ Br = $true },\n  @{ Text = '## Display:' },\n  @{ Br = $true },\n  @{ Text = '- 6.7-inch Super Retina XDR display' },\n  @{ Br = $true },\n  @{ Text = '- Always-On and ProMotion technology with adaptive 120Hz refresh rate' },\n  @{ Br = $true },\n  @{ Text = '- Oleophobic coating' },\n  @{ Br = $true },\n  @{ Text = '- 2796x1290 resolution at 460 ppi' },\n  @{ Br = $true },\n  @{ Text = '- HDR support' },\n  @{ Br = $true },\n  @{ Text = '## ' },\n  @{ Br = $true },\n  @{ Text = '## Camera:' },\n  @{ Br = $true },\n  @{ Text = '- 48MP Main camera with 24mm \u0192/1.78 aperture' },\n  @{ Br = $true },\n  @{ Text = '- 12MP Ultra Wide 13mm \u0192/2.2 aperture camera' },\n  @{ Br = $true },\n  @{ Text = '- 12MP 2x Telephoto 48mm \u0192/1.78 aperture camera enabled by quad-pixel sensor' },\n  @{ Br = $true },\n  @{ Text = '- 12MP 3x Telephoto 77mm \u0192/2.8 aperture camera' },\n  @{ Br = $true },\n  @{ Text = '- Night mode portraits enabled by LiDAR Scanner' },\n  @{ Br = $true },\n  @{ Text = '- Portrait mode with advanced bokeh and Depth Control' },\n  @{ Br = $true },\n  @{ Text = '- 4K video recording at 24, 25, 30 or 60 fps' },\n  @{ Br = $true },\n  @{ Text = '## ' },\n  @{ Br = $true },\n  @{ Text = '## Body:' },\n  @{ Br = $true },\n  @{ Text = '- Ceramic Shield front, textured matte glass back and stainless steel design' },\n  @{ Br = $true },\n  @{ Text = '- Resists splashes, water and dust with IP68 rating' },\n  @{ Br = $true },\n  @{ Text = '- Available in Space Black, Silver, Gold and Deep Purple' },\n  @{ Br = $true },\n  @{ Text = '- Dimensions: 3.05 inches x 6.33 inches x 0.31 inch (77.6 mm x 160.7 mm x 7.85 mm)' },\n  @{ Br = $true },\n  @{ Text = '- Weight: 8.47 ounces (240 grams)' },\n  @{ Br = $true },\n  @{ Text = '## ' },\n  @{ Br = $true },\n  @{ Text = '## Power:' },\n  @{ Br = $true },\n  @{ Text = '- Built-in rechargeable lithium-ion battery with up to 29 hours video playback' },\n  @{ Br = $true },\n  @{ Text = '- MagSafe and Qi wireless charging up to 15W' },\n  @{ Br = $true },\n  @{ Text = '- Fast charging capable: up to 50% in 30 minutes with 20W adapter or higher' },\n  @{ Br = $true },\n  @{ Text = '## ' },\n  @{ Br = $true },\n  @{ Text = '## Storage:' },\n  @{ Br = $true },\n  @{ Text = '- 128GB, 256GB, 512GB and 1TB versions available' },\n  @{ Br = $true },\n  @{ Text = '## ' },\n  @{ Br = $true },\n  @{ Text = '## Connectivity:' },\n  @{ Br = $true },\n  @{ Text = '- 5G, 4G LTE, Wi-Fi 6, Bluetooth\u00ae and NFC' },\n  @{ Br = $true },\n  @{ Text = '- Dual eSIM support, no physical SIM slot' },\n  @{ Br = $true },\n  @{ Text = '- GPS, GLONASS, Galileo, QZSS and BeiDou' },\n  @{ Br = $true },\n  @{ Text = '## ' },\n  @{ Br = $true },\n  @{ Text = '## Price:' },\n  @{ Br = $true },\n  @{ Text = '- $1,099 (128GB), $1,199 (256GB), $1,399 (512GB), $1,599 (1TB)' }\n)\n\nfunction ConvertTo-XmlText {\n    param([string]$Value)\n    $escaped = $Value.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')\n    if ($Value -ne $Value.Trim() -or $Value.Length -eq 0) {\n        return '<w:t xml:space=\"preserve\">' + $escaped + '</w:t>'\n    } else {\n        return '<w:t>' + $escaped + '</w:t>'\n    }\n}\n\n$sb = New-Object System.Text.StringBuilder\nforeach ($seg in $segments) {\n    if ($seg.Br) {\n        [void]$sb.Append('<w:br/>')\n    } else {\n        [void]$sb.Append((ConvertTo-XmlText $seg.Text))\n    }\n}\n$runInner = $sb.ToString()\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r>' + $runInner + '</w:r></w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$r = $d.Content\n$r.InsertXML($ooxml)\n"}
